$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells stay text-typed (matching original inlineStr cells) before assigning values

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.583.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.580.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.08%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.80"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.034.90"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.43"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.509.19"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.582.31"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.494"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.77"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.05"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0838"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.92"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.05"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.30"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.18%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.74"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.19"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.67"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.80%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "336.26"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.03"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.892"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.93"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.39"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.47"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.147.73"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.606"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.92"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0545"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.40%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.55"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0967"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0239"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.31%  "
